$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Devices")

# Update the "Other Slot Cards" count text in N10 and N11 from (0 of 5) to (0 of 17)
$ws.Range("N10").Value = "Other Slot Cards  (0 of 17)"
$ws.Range("N11").Value = "Other Slot Cards  (0 of 17)"

# Update the view's selection to match the new state
$ws.Range("N11").Select()
